$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.550.98'
$ws.Range("E2").Value = '  +2.19%  '
$ws.Range("D3").Value = '1.674.41'
$ws.Range("E3").Value = '  +1.88%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'220.61"
$ws.Range("E5").Value = '  +2.66%  '
$ws.Range("D6").Value = "'0.5287"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = "'0.2684"
$ws.Range("E8").Value = '  +3.00%  '
$ws.Range("D9").Value = "'0.06394"
$ws.Range("E9").Value = '  +0.58%  '
$ws.Range("D10").Value = "'21.92"
$ws.Range("E10").Value = '  +5.84%  '
$ws.Range("D11").Value = "'0.07803"
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("D12").Value = "'4.494"
$ws.Range("E12").Value = '  +1.58%  '
$ws.Range("D13").Value = '1.669.72'
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("D14").Value = "'0.5583"
$ws.Range("E14").Value = '  +1.09%  '
$ws.Range("D15").Value = '0.0₅8337'
$ws.Range("E15").Value = '  +0.74%  '
$ws.Range("D16").Value = "'65.82"
$ws.Range("E16").Value = '  +1.72%  '
$ws.Range("D17").Value = '26.536.27'
$ws.Range("E17").Value = '  +2.14%  '
$ws.Range("D19").Value = "'4.775"
$ws.Range("E19").Value = '  +1.60%  '
$ws.Range("D20").Value = "'193.58"
$ws.Range("E20").Value = '  +2.33%  '
$ws.Range("D21").Value = "'10.37"
$ws.Range("E21").Value = '  +1.97%  '
$ws.Range("D22").Value = "'6.321"
$ws.Range("E22").Value = '  +0.98%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = "'0.1272"
$ws.Range("E24").Value = '  +2.41%  '
$ws.Range("D25").Value = "'139.78"
$ws.Range("E25").Value = '  -3.10%  '
$ws.Range("D26").Value = "'7.417"
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("D27").Value = "'16.34"
$ws.Range("E27").Value = '  +2.82%  '
$ws.Range("D28").Value = "'1.425"
$ws.Range("E28").Value = '  +2.43%  '
$ws.Range("D29").Value = "'0.06211"
$ws.Range("E29").Value = '  +5.00%  '
$ws.Range("E30").Value = '  +2.47%  '
$ws.Range("D31").Value = "'3.622"
$ws.Range("E31").Value = '  +6.74%  '
$ws.Range("D32").Value = "'3.434"
$ws.Range("E32").Value = '  +1.26%  '
$ws.Range("D33").Value = "'1.686"
$ws.Range("E33").Value = '  +2.37%  '
$ws.Range("D34").Value = "'1.011"
$ws.Range("E34").Value = '  +1.79%  '
$ws.Range("D35").Value = "'0.6094"
$ws.Range("E35").Value = '  +8.32%  '
$ws.Range("D36").Value = "'2.415"
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("D37").Value = "'2.780"
$ws.Range("E37").Value = '  +1.25%  '
$ws.Range("E38").Value = '  +1.00%  '
$ws.Range("D39").Value = "'6.070"
$ws.Range("E39").Value = '  +3.46%  '
$ws.Range("D40").Value = '1.094.09'
$ws.Range("E40").Value = '  +5.98%  '
$ws.Range("D41").Value = "'0.8582"
$ws.Range("E41").Value = '  +0.75%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").Value = "'100.73"
$ws.Range("E43").Value = '  +1.98%  '
$ws.Range("D44").Value = '1.818.63'
$ws.Range("E44").Value = '  +1.52%  '
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("D46").Value = "'58.54"
$ws.Range("E46").Value = '  +5.25%  '
$ws.Range("E47").Value = '  +0.44%  '
$ws.Range("D48").Value = "'8.112"
$ws.Range("E48").Value = '  +0.93%  '
$ws.Range("D49").Value = "'1.517"
$ws.Range("E49").Value = '  +10.70%  '
$ws.Range("D50").Value = "'0.05199"
$ws.Range("E50").Value = '  +1.04%  '
$ws.Range("D51").Value = "'6.020"
$ws.Range("E51").Value = '  +2.17%  '
